# pt_element: instead of composing ec class for each span in a term,
# compose it once (as class_nl, formerly pre_class) and use it repeatedly.
#
# Add a new worksheet "ec_class accumulated" in front of the existing
# sheets, holding the accumulated-timing comparison between the old
# "construct class for each term" approach and the new
# "construct pre_class once" approach.

$wb = $excel.ActiveWorkbook

$first = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($first)
$ws.Name = "ec_class accumulated"

# Column widths match the other per-run comparison sheet ("db query vs get")
$ws.Columns.Item(2).ColumnWidth = 23
$ws.Columns.Item(3).ColumnWidth = 3.57

# Headers
$ws.Range("B1").Value = "pt_element: construct class for each term"
$ws.Range("D1").Value = "pt_element: construct pre_class once"

$ws.Range("A2").Value = "Run"
$ws.Range("B2").Value = "Total pt_element_accumulator"
$ws.Range("D2").Value = "Total pt_element_accumulator"

$ws.Range("A3").Value = "Average"
$ws.Range("B3").Formula = "=AVERAGE(B6:B32)"
$ws.Range("D3").Formula = "=AVERAGE(D6:D32)"

$ws.Range("A4").Value = "StDev"
$ws.Range("B4").Formula = "=STDEVA(B6:B32)"
$ws.Range("D4").Formula = "=STDEVA(D6:D32)"

$ws.Range("A5").Value = "RelStdDev"
$ws.Range("B5").Formula = "=B4/B3"
$ws.Range("D5").Formula = "=D4/D3"
$ws.Range("B5").NumberFormat = "0.00%"
$ws.Range("D5").NumberFormat = "0.00%"

# Raw per-run samples
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0.42475299999999999
$ws.Range("D6").Value = 0.39991199999999999

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 0.438661999999999
$ws.Range("D7").Value = 0.45086799999999899

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 0.446709999999999
$ws.Range("D8").Value = 0.484957999999999

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 0.42122699999999902
$ws.Range("D9").Value = 0.44270900000000002

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 0.47458699999999998
$ws.Range("D10").Value = 0.47355799999999898

$ws.Range("B31").Select() | Out-Null
